$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.213
$ws.Range("B4").Value = 4.627800000000001
$ws.Range("B5").Value = 5.169299999999998
$ws.Range("A7").Value = -21.62710000000001
$ws.Range("B8").Value = 5.0171
$ws.Range("A16").Value = -21.52870000000001
$ws.Range("B16").Value = 4.749300000000003
